$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates ---
$ws.Range("D9").Value = "1399-03-12 (9)"
$ws.Range("E9").Value = "1400-02-29 (8)"
$ws.Range("F9").Value = "1401-02-27 (11)"
$ws.Range("G9").Value = "1402-02-28 (7)"
$ws.Range("H9").Value = "1402-02-28"

# --- Data rows 11-27 ---
# Row 11
$ws.Range("D11").Value = 20582
$ws.Range("E11").Value = 15758
$ws.Range("F11").Value = 14557
$ws.Range("G11").Value = 19296
$ws.Range("H11").Value = 20800

# Row 12
$ws.Range("D12").Value = -13328
$ws.Range("E12").Value = -11514
$ws.Range("F12").Value = -9416
$ws.Range("G12").Value = -13316
$ws.Range("H12").Value = -16725

# Row 13
$ws.Range("D13").Value = 7254
$ws.Range("E13").Value = 4245
$ws.Range("F13").Value = 5140
$ws.Range("G13").Value = 5980
$ws.Range("H13").Value = 4075

# Row 14
$ws.Range("D14").Value = -1418
$ws.Range("E14").Value = -1378
$ws.Range("F14").Value = -1116
$ws.Range("G14").Value = -1653
$ws.Range("H14").Value = -2009

# Row 15
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

# Row 16
$ws.Range("D16").Value = -41
$ws.Range("E16").Value = -126
$ws.Range("F16").Value = -179
$ws.Range("G16").Value = -176
$ws.Range("H16").Value = -294

# Row 17
$ws.Range("D17").Value = 5795
$ws.Range("E17").Value = 2741
$ws.Range("F17").Value = 3846
$ws.Range("G17").Value = 4150
$ws.Range("H17").Value = 1773

# Row 18
$ws.Range("D18").Value = -2383
$ws.Range("E18").Value = -2093
$ws.Range("F18").Value = -1436
$ws.Range("G18").Value = -1933
$ws.Range("H18").Value = -1965

# Row 19
$ws.Range("D19").Value = 4981
$ws.Range("E19").Value = 8303
$ws.Range("F19").Value = 5332
$ws.Range("G19").Value = 7921
$ws.Range("H19").Value = 11203

# Row 20
$ws.Range("D20").Value = 8393
$ws.Range("E20").Value = 8951
$ws.Range("F20").Value = 7743
$ws.Range("G20").Value = 10138
$ws.Range("H20").Value = 11010

# Row 21
$ws.Range("D21").Value = -852
$ws.Range("E21").Value = -370
$ws.Range("F21").Value = -653
$ws.Range("G21").Value = -536
$ws.Range("H21").Value = "-"

# Row 22
$ws.Range("D22").Value = 7541
$ws.Range("E22").Value = 8581
$ws.Range("F22").Value = 7090
$ws.Range("G22").Value = 9602
$ws.Range("H22").Value = 11010

# Row 23
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 5

# Row 24
$ws.Range("D24").Value = 7541
$ws.Range("E24").Value = 8581
$ws.Range("F24").Value = 7090
$ws.Range("G24").Value = 9606
$ws.Range("H24").Value = 11015

# Row 25
$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

# Row 26
$ws.Range("D26").Value = "-"
$ws.Range("E26").Value = 8840
$ws.Range("F26").Value = 10031
$ws.Range("G26").Value = 8596
$ws.Range("H26").Value = 6427

# Row 27
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0

